$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C become C/D)
$ws.Columns.Item(2).Insert()

# New column B header and width
$ws.Range("B1").Value = "StatQuery"
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()

# Fill in the two query cells in row 2 (A2 was previously blank w/ wrap style, B2 is new)
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_type IN [''Variants file'']  RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("B2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_type IN [''Aligned DNA reads file'',''Aligned RNA reads file'',''Index file'',''Variants file''] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trials'
$ws.Range("A2:B2").WrapText = $true

# Row height for the query row
$ws.Rows.Item(2).RowHeight = 101.5

# Update selection/view to match target (select A2, which also clears topLeftCell scroll)
$ws.Range("A2").Select()
